$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1851
$ws.Range("F3").Value = 407
$ws.Range("F4").Value = 1503
$ws.Range("F6").Value = 388
$ws.Range("F7").Value = 751
$ws.Range("F8").Value = 13245
$ws.Range("F9").Value = 13132
$ws.Range("F10").Value = 1004
$ws.Range("F13").Value = 551
$ws.Range("F15").Value = 653
$ws.Range("F16").Value = 2078
$ws.Range("F17").Value = 60
$ws.Range("F19").Value = 54
$ws.Range("F21").Value = 218
$ws.Range("F22").Value = 278
$ws.Range("F23").Value = 746

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 59
$ws.Range("F9").Value = 19

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 22

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1851
$ws.Range("F4").Value = 407
$ws.Range("F5").Value = 1503
$ws.Range("F7").Value = 388
$ws.Range("F9").Value = 751
$ws.Range("F10").Value = 13245
$ws.Range("F11").Value = 13132
$ws.Range("F12").Value = 1004
$ws.Range("F15").Value = 551
$ws.Range("F17").Value = 653
$ws.Range("F20").Value = 2078
$ws.Range("F21").Value = 60
$ws.Range("F23").Value = 54
$ws.Range("F26").Value = 59
$ws.Range("F27").Value = 22
$ws.Range("F28").Value = 218
$ws.Range("F29").Value = 278
$ws.Range("F30").Value = 746
$ws.Range("F33").Value = 19

$wb.Save()
